# Applies the "my latest familt tree Excel file" edit described by the diff:
#  - workbook.xml: workbookPr/workbookView cosmetic attributes (best effort; these
#    are window/view state that the host app may not round-trip through the xlsx
#    writer, but we still set them through the object model where exposed)
#  - sharedStrings/worksheet: 58 cell text updates on "Family Tree - All" (mostly
#    appending "b: #### d: ####"-style birth/death placeholders, or full name/date
#    rewrites) plus a new value in a previously-empty cell (F167), and the saved
#    selection/scroll position.

$wb = $excel.ActiveWorkbook

# --- Workbook-level cosmetic properties -----------------------------------------
try { $wb.CheckCompatibility = $true } catch { }
try {
    $win = $excel.ActiveWindow
    $win.Left = 0
    $win.Top = 0
    $win.Width = 1530
    $win.Height = 803
} catch { }

# --- Worksheet cell updates on "Family Tree - All" -------------------------------
$ws = $wb.Worksheets.Item("Family Tree - All")

$ws.Range("D95").Value = "Kenneth the Alpin (####-####)"
$ws.Range("D96").Value = "Constantine (####-####)"
$ws.Range("D97").Value = "Donald IV (####-####)"
$ws.Range("D98").Value = "Malcolm (####-####)"
$ws.Range("D99").Value = "Kenneth II (####-####)"
$ws.Range("D100").Value = "Malcolm II (####-####)"
$ws.Range("D101").Value = "Beatrix (####-####)"
$ws.Range("D102").Value = "Duncan (####-####)"
$ws.Range("D103").Value = "Malcolm (####-####)"
$ws.Range("D104").Value = "David (####-####)"
$ws.Range("D105").Value = "Henry (####-####)"
$ws.Range("D106").Value = "David (####-####)"
$ws.Range("D107").Value = "Isabel (####-####)"
$ws.Range("D108").Value = "Robert (####-####)"
$ws.Range("D109").Value = "Robert (####-####)"
$ws.Range("D110").Value = "Robert the Bruce (####-####)"
$ws.Range("D111").Value = "Marlory (####-####)"
$ws.Range("D112").Value = "Robert II (####-####)"
$ws.Range("D113").Value = "Robert III (####-####)"
$ws.Range("D114").Value = "James I (####-####)"
$ws.Range("D115").Value = "James II (####-####)"
$ws.Range("D116").Value = "James III (####-####)"
$ws.Range("D117").Value = "James V (####-####)"
$ws.Range("D118").Value = "Mary Queen of Scots (####-####)"
$ws.Range("D119").Value = "James VI (####-####)"
$ws.Range("D120").Value = "Elizabeth Stewart (####-####)"
$ws.Range("D121").Value = "Sophia (####-####)"
$ws.Range("D122").Value = "George I (####-####)"
$ws.Range("D123").Value = "George II (####-####)"
$ws.Range("D124").Value = "Frederick Prince of Wales (####-####)"
$ws.Range("D125").Value = "George III (####-####)"
$ws.Range("D126").Value = "Edward Duke of Kent (####-####)"
$ws.Range("D127").Value = "George IV (####-####)"
$ws.Range("D128").Value = "George V (####-####)"
$ws.Range("E128").Value = "Mary II (####-####)"
$ws.Range("E150").Value = "Thomas Holland (####-####)"
$ws.Range("E159").Value = "Joan Or Jane De VALLETORT b: Abt 1213 d: ####"
$ws.Range("D160").Value = "Lawrence De CORNWALL b: Abt 1241 d: ####"
$ws.Range("D161").Value = "Sybil DE CORNWALL b: #### d: ####"
$ws.Range("E163").Value = "Lady Maud DE PENNINGTON b: #### d: ####"
$ws.Range("E165").Value = "Anne Lady of C M Fenwicke b: 1403 d: ####              OR"
$ws.Range("E166").Value = "Joan UNKNOWN (Stapleton?) b: #### d: ####           OR"
$ws.Range("E167").Value = "Mary FENWICK b: #### d: ####                                     OR "
$ws.Range("F167").Value = "Joan Stapleton (?) b: #### d: ####"
$ws.Range("D169").Value = "John HUDDLESTONE b: 1490 d: ####"
$ws.Range("E169").Value = "Elizabeth SUTTON b: 1493 d: ####"
$ws.Range("E170").Value = "Bridget COTTON b: 1530 d: ####"
$ws.Range("E171").Value = "Dorthy BEACONSALL b: 1552 d: ####"
$ws.Range("D172").Value = "Henry HUDDLESTON b: 1575 d: ####"
$ws.Range("E172").Value = "Dorthy DORMER b: 1577 d: ####"
$ws.Range("D174").Value = "Richard HUDDLESTON b: Abt 1660 d: ####"
$ws.Range("E174").Value = "Mary BOSTOCK b: 1660 d: ####"
$ws.Range("E177").Value = "Mary PATTERSON b: #### d: ####"
$ws.Range("E178").Value = "Sarah GALLAWAY b: #### d: ####"
$ws.Range("D179").Value = "William Thomas (1793-1834)"
$ws.Range("E179").Value = "Nancy Huddleston (1794-1865)"
$ws.Range("D185").Value = "Gregory Alan (Johnson) Jonason (1973-####)"
$ws.Range("E185").Value = "Marie Kenly (Antoine) Jonason (1972-####)"

# --- Restore the saved selection / scroll position --------------------------------
$ws.Activate()
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 149
    $win.ScrollColumn = 1
} catch { }
$ws.Range("D179").Select()

